# Refresh the crypto price/volume table (GitHub Actions data pull).
# D-column numeric-looking strings are written with a leading "'"
# (Excel quote-prefix) so they stay literal text -- exactly like the
# original sheet -- instead of silently being parsed into numbers
# and losing formatting such as trailing zeros ("64.00" -> 64).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.131.77"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "2.575.80"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "'561.69"
$ws.Range("E5").Value = "  -0.85%  "
$ws.Range("D6").Value = "'142.60"
$ws.Range("E6").Value = "  -2.16%  "
$ws.Range("E7").Value = "  +0.28%  "
$ws.Range("D8").Value = "'0.597"
$ws.Range("E8").Value = "  -1.11%  "
$ws.Range("D9").Value = "2.578.12"
$ws.Range("E9").Value = "  -2.71%  "
$ws.Range("D10").Value = "'6.64"
$ws.Range("E10").Value = "  -2.64%  "
$ws.Range("E11").Value = "  -0.34%  "
$ws.Range("D12").Value = "'0.159"
$ws.Range("E12").Value = "  +11.73%  "
$ws.Range("D13").Value = "'0.344"
$ws.Range("E13").Value = "  +0.81%  "
$ws.Range("D14").Value = "3.027.34"
$ws.Range("E14").Value = "  -2.30%  "
$ws.Range("D15").Value = "59.118.03"
$ws.Range("E15").Value = "  -1.63%  "
$ws.Range("D16").Value = "'22.49"
$ws.Range("E16").Value = "  +3.77%  "
$ws.Range("E17").Value = "  +0.39%  "
$ws.Range("D18").Value = "2.572.51"
$ws.Range("E18").Value = "  -2.96%  "
$ws.Range("D19").Value = "'4.53"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "'336.52"
$ws.Range("E20").Value = "  -2.04%  "
$ws.Range("D21").Value = "'10.31"
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").Value = "'6.32"
$ws.Range("E22").Value = "  +0.76%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'64.00"
$ws.Range("E24").Value = "  -4.17%  "
$ws.Range("D25").Value = "'0.457"
$ws.Range("E25").Value = "  +5.12%  "
$ws.Range("E26").Value = "  +0.44%  "
$ws.Range("E27").Value = "  -2.46%  "
$ws.Range("D28").Value = "'7.26"
$ws.Range("E28").Value = "  -0.61%  "
$ws.Range("D29").Value = "0.0₃0774"
$ws.Range("E29").Value = "  +0.33%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  +0.11%  "

# Rows 31/32 exchange rank order (Monero moves above PancakeSwap).
# Column A rank numbers are untouched; only B:E are rewritten with
# each coin's own refreshed data, in its new row.
$ws.Range("B31").Value = "Monero"
$ws.Range("C31").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D31").Value = "'161.18"
$ws.Range("E31").Value = "  +2.05%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "'1.67"
$ws.Range("E32").Value = "  -2.19%  "

$ws.Range("D33").Value = "'6.06"
$ws.Range("E33").Value = "  -0.80%  "
$ws.Range("E34").Value = "  -1.07%  "
$ws.Range("D35").Value = "'4.00"
$ws.Range("E35").Value = "  -1.67%  "
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "'0.880"
$ws.Range("E37").Value = "  -3.93%  "
$ws.Range("D38").Value = "'37.42"
$ws.Range("E38").Value = "  -0.15%  "
$ws.Range("D39").Value = "'0.869"
$ws.Range("E39").Value = "  -4.68%  "
$ws.Range("E40").Value = "  -2.72%  "
$ws.Range("D41").Value = "'293.86"
$ws.Range("E41").Value = "  -2.27%  "
$ws.Range("D42").Value = "'3.64"
$ws.Range("E42").Value = "  -0.49%  "
$ws.Range("D43").Value = "'0.998"
$ws.Range("E43").Value = "  +0.30%  "
$ws.Range("D44").Value = "'131.44"
$ws.Range("E44").Value = "  +9.15%  "
$ws.Range("D45").Value = "'0.0972"
$ws.Range("E45").Value = "  -0.15%  "
$ws.Range("D46").Value = "'0.593"
$ws.Range("E46").Value = "  -1.81%  "
$ws.Range("D47").Value = "'10.63"
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("E48").Value = "  -2.09%  "
$ws.Range("D49").Value = "'19.03"
$ws.Range("E49").Value = "  -1.83%  "
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("D51").Value = "'18.38"
$ws.Range("E51").Value = "  -0.38%  "
